$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B13").Value = "Tamby"
$ws.Range("B14").Value = "Tamby"

$ws.Range("B15").Select()
